# Swap the order of the recorder names in the "Recorded By" column (G).
# Every cell that currently reads "dnasr281@gmail.com, System" should be
# changed to "System, dnasr281@gmail.com". All other values in that column
# (e.g. "System", "dnasr281@gmail.com", or "System, dnasr281@gmail.com")
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
